$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 529; this shifts rows 529:644 down to 530:645
$ws.Rows("529:529").Insert()

# Populate the freshly inserted row 529 with the new record
$ws.Cells.Item(529, 1).Value = 9
$ws.Cells.Item(529, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(529, 3).Value = "Metropolitana"
$ws.Cells.Item(529, 4).Value = 45244
$ws.Cells.Item(529, 5).Value = 13
$ws.Cells.Item(529, 6).Value = 100112044
$ws.Cells.Item(529, 7).Value = "Perejil"
$ws.Cells.Item(529, 8).Value = "Sin especificar"
$ws.Cells.Item(529, 9).Value = "Primera"
$ws.Cells.Item(529, 10).Value = 70
$ws.Cells.Item(529, 11).Value = 12000
$ws.Cells.Item(529, 12).Value = 13000
$ws.Cells.Item(529, 13).Value = 12500
$ws.Cells.Item(529, 14).Value = "`$/docena de atados"
$ws.Cells.Item(529, 15).Value = "Región Metropolitana"
$ws.Cells.Item(529, 16).Value = 4167
$ws.Cells.Item(529, 17).Value = 3
$ws.Cells.Item(529, 18).Value = "Hortaliza"
